$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (header + 16 players), reordered/trimmed per the updated
# "Team of Outs" list (Bobby Portis and Isaiah Stewart removed).
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Vasilije Micic", "PG,SG", "Charlotte Hornets"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The sheet previously had 18 rows of data (19 rows total incl. header).
# The updated roster only has 16 rows of data (17 total), so remove the
# two now-unused trailing rows.
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()
